$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Summary")
$ws.Range("B4").Value = "inf"
$ws.Range("B6").Value = -56101.34517283511
$ws.Range("B7").Value = 9909345.752283832
$ws.Range("B8").Value = 24184978.37650956
$ws.Range("B10").Value = 3451084.313349553

$ws = $wb.Worksheets.Item("Costs and Revenues")
$ws.Range("E2").Value = 47447.79405033727
$ws.Range("F2").Value = 47665.20296074481
$ws.Range("G2").Value = 53820.30179896292
$ws.Range("H2").Value = 55588.10794580438
$ws.Range("I2").Value = 40872.6284777889
$ws.Range("J2").Value = 40068.39727357774
$ws.Range("K2").Value = 43802.19989414858
$ws.Range("L2").Value = 39848.08030419912
$ws.Range("M2").Value = 39146.38296673041
$ws.Range("N2").Value = 30334.34295949766
$ws.Range("O2").Value = 33415.60698489296
$ws.Range("P2").Value = 27138.64542625533
$ws.Range("E3").Value = 133100.0000000001
$ws.Range("E4").Value = 22142.23997720535
$ws.Range("F4").Value = 22359.64888761288
$ws.Range("G4").Value = 28514.747725831
$ws.Range("H4").Value = 30282.55387267243
$ws.Range("I4").Value = 15567.07440465696
$ws.Range("J4").Value = 14762.84320044582
$ws.Range("K4").Value = 18496.64582101665
$ws.Range("L4").Value = 14542.52623106719
$ws.Range("M4").Value = 13840.82889359847
$ws.Range("N4").Value = 5028.788886365757
$ws.Range("O4").Value = 8110.052911761053
$ws.Range("P4").Value = 1833.091353123424
$ws.Range("F6").Value = 21693.9419681236
$ws.Range("H6").Value = 21693.94196812362
$ws.Range("I6").Value = 21693.94196812361
$ws.Range("K6").Value = 21693.94196812361
$ws.Range("L6").Value = 21693.94196812361
$ws.Range("M6").Value = 21693.94196812361
$ws.Range("N6").Value = 21693.94196812358
$ws.Range("O6").Value = 21693.94196812358
$ws.Range("P6").Value = 21693.94196812358

$ws = $wb.Worksheets.Item("Fed-in Capacity")
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = 0
$ws.Range("O11").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = 0
$ws.Range("O12").Value = 0
$ws.Range("P12").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = 0
$ws.Range("O14").Value = 0
$ws.Range("L15").Value = 61.18167021676314
$ws.Range("L17").Value = 130.6648563030561
$ws.Range("M17").Value = 113.4004983079896
$ws.Range("L18").Value = 61.18167021676314
$ws.Range("M18").Value = 51.84373129681028
$ws.Range("N18").Value = 38.66169381481656
$ws.Range("O18").Value = 57.81213424001893
$ws.Range("P18").Value = 65.92768427608706
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = 0
$ws.Range("N19").Value = 0
$ws.Range("O19").Value = 0
$ws.Range("P19").Value = 0
$ws.Range("M20").Value = 113.4004983079896
$ws.Range("N20").Value = 110.5750244233121
$ws.Range("L21").Value = 61.18167021676314
$ws.Range("L23").Value = 130.6648563030561
$ws.Range("N23").Value = 110.5750244233121
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 80.29914934735042
$ws.Range("L24").Value = 61.18167021676314
$ws.Range("M24").Value = 51.84373129681028
$ws.Range("J25").Value = 0
$ws.Range("Q25").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("N27").Value = 38.66169381481656
$ws.Range("O27").Value = 57.81213424001893
$ws.Range("P27").Value = 0
$ws.Range("Q27").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 61.18167021676314
$ws.Range("M30").Value = 51.84373129681028
$ws.Range("O30").Value = 57.81213424001893
$ws.Range("M32").Value = 0
$ws.Range("K33").Value = 80.29914934735042
$ws.Range("P33").Value = 65.92768427608706
$ws.Range("Q33").Value = 94.49434172313325
$ws.Range("K36").Value = 0
$ws.Range("M36").Value = 51.84373129681028
$ws.Range("K39").Value = 0
$ws.Range("Q39").Value = 0
$ws.Range("Q41").Value = 0
$ws.Range("J42").Value = 93.17061249236157
$ws.Range("Q44").Value = 0
$ws.Range("R44").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 80.29914934735042
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 0

$ws = $wb.Worksheets.Item("Unmet Demand")
$ws.Range("L11").Value = 130.6648563030561
$ws.Range("M11").Value = 113.4004983079896
$ws.Range("O11").Value = 117.8828208804077
$ws.Range("L12").Value = 61.18167021676314
$ws.Range("M12").Value = 51.84373129681028
$ws.Range("O12").Value = 57.81213424001893
$ws.Range("P12").Value = 65.92768427608706
$ws.Range("L14").Value = 130.6648563030561
$ws.Range("M14").Value = 113.4004983079896
$ws.Range("O14").Value = 117.8828208804077
$ws.Range("L15").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = 0
$ws.Range("N18").Value = 0
$ws.Range("O18").Value = 0
$ws.Range("P18").Value = 0
$ws.Range("K19").Value = 94.30397654773019
$ws.Range("L19").Value = 90.4687457914608
$ws.Range("M19").Value = 92.09541281912071
$ws.Range("N19").Value = 81.96869489115805
$ws.Range("O19").Value = 96.22962838366004
$ws.Range("P19").Value = 101.5955875616828
$ws.Range("M20").Value = 0
$ws.Range("N20").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").Value = 0
$ws.Range("J24").Value = 93.17061249236157
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = 0
$ws.Range("J25").Value = 105.873818686614
$ws.Range("Q25").Value = 126.4887893424616
$ws.Range("K26").Value = 135.370731907559
$ws.Range("J27").Value = 93.17061249236157
$ws.Range("N27").Value = 0
$ws.Range("O27").Value = 0
$ws.Range("P27").Value = 65.92768427608706
$ws.Range("Q27").Value = 94.49434172313325
$ws.Range("K28").Value = 94.30397654773019
$ws.Range("L28").Value = 90.4687457914608
$ws.Range("K29").Value = 135.370731907559
$ws.Range("L29").Value = 130.6648563030561
$ws.Range("K30").Value = 80.29914934735042
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = 0
$ws.Range("O30").Value = 0
$ws.Range("M32").Value = 113.4004983079896
$ws.Range("K33").Value = 0
$ws.Range("P33").Value = 0
$ws.Range("Q33").Value = 0
$ws.Range("K36").Value = 80.29914934735042
$ws.Range("M36").Value = 0
$ws.Range("K39").Value = 80.29914934735042
$ws.Range("Q39").Value = 94.49434172313325
$ws.Range("Q41").Value = 150.3839754851235
$ws.Range("J42").Value = 0
$ws.Range("Q44").Value = 150.3839754851235
$ws.Range("R44").Value = 173.7492132756177
$ws.Range("J45").Value = 93.17061249236157
$ws.Range("K45").Value = 0
$ws.Range("J46").Value = 105.873818686614
$ws.Range("K46").Value = 94.30397654773019

$ws = $wb.Worksheets.Item("Household Surplus")
$ws.Range("B5").Value = 306829.0042792289
$ws.Range("B6").Value = 308350.8666520815
$ws.Range("B7").Value = 351436.5585196086
$ws.Range("B8").Value = 363811.2015474986
$ws.Range("B9").Value = 260802.84527139
$ws.Range("B10").Value = 255173.226841912
$ws.Range("B11").Value = 281309.8451859078
$ws.Range("B12").Value = 253631.0080562617
$ws.Range("B13").Value = 248719.1266939807
$ws.Range("B14").Value = 187034.8466433517
$ws.Range("B15").Value = 208603.6948211186
$ws.Range("B16").Value = 164664.9639106554
